$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.814.28"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").Value = "3.786.61"
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'443.47"
$ws.Range("E5").Value = "  +5.14%  "

$ws.Range("D6").Value = "'143.86"
$ws.Range("E6").Value = "  +12.76%  "

$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  +2.77%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.730"
$ws.Range("E9").Value = "  +2.04%  "

$ws.Range("E10").Value = "  -8.14%  "

$ws.Range("E11").Value = "  -10.26%  "

$ws.Range("D12").Value = "'43.15"
$ws.Range("E12").Value = "  +7.95%  "

$ws.Range("D13").Value = "'10.31"
$ws.Range("E13").Value = "  +4.32%  "

$ws.Range("D14").Value = "4.388.05"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").Value = "'14.69"
$ws.Range("E15").Value = "  -6.56%  "

$ws.Range("D16").Value = "3.818.61"
$ws.Range("E16").Value = "  -0.37%  "

$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Value = "'19.78"
$ws.Range("E18").Value = "  +1.82%  "

$ws.Range("E19").Value = "  +6.99%  "

$ws.Range("D20").Value = "66.883.11"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("D21").Value = "'412.47"
$ws.Range("E21").Value = "  +2.90%  "

$ws.Range("D22").Value = "'14.48"
$ws.Range("E22").Value = "  +2.26%  "

$ws.Range("E23").Value = "  +10.03%  "

$ws.Range("D24").Value = "'85.27"
$ws.Range("E24").Value = "  +2.16%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'3.40"
$ws.Range("E25").Value = "  +7.18%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'36.82"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "'9.60"
$ws.Range("E27").Value = "  +29.13%  "

$ws.Range("D28").Value = "'5.50"
$ws.Range("E28").Value = "  -4.57%  "

$ws.Range("D29").Value = "'9.68"
$ws.Range("E29").Value = "  +3.57%  "

$ws.Range("D30").Value = "'732.37"
$ws.Range("E30").Value = "  +4.27%  "

$ws.Range("D31").Value = "'13.83"
$ws.Range("E31").Value = "  +13.27%  "

$ws.Range("D32").Value = "'0.134"
$ws.Range("E32").Value = "  +11.32%  "

$ws.Range("D33").Value = "'2.70"
$ws.Range("E33").Value = "  -1.73%  "

$ws.Range("D34").Value = "'43.33"
$ws.Range("E34").Value = "  +15.35%  "

$ws.Range("E35").Value = "  +6.68%  "

$ws.Range("D36").Value = "'56.42"
$ws.Range("E36").Value = "  +3.32%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").Value = "'5.48"
$ws.Range("E38").Value = "  +24.90%  "

$ws.Range("D39").Value = "'0.0475"
$ws.Range("E39").Value = "  +5.93%  "

$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "'2.88"
$ws.Range("E40").Value = "  -0.76%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = "  +28.99%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.140"
$ws.Range("E42").Value = "  +4.83%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.336"
$ws.Range("E44").Value = "  +17.74%  "

$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0668"
$ws.Range("E45").Value = "  -11.43%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.29"
$ws.Range("E46").Value = "  +6.65%  "

$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "  +0.75%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'144.64"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'2.08"
$ws.Range("E49").Value = "  +1.64%  "

$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "'2.64"
$ws.Range("E50").Value = "  +5.15%  "

$ws.Range("E51").Value = "  +4.14%  "
